$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192543387413025
$ws.Range("B1").Value = 3.401607036590576
$ws.Range("C1").Value = 4.34290075302124
$ws.Range("D1").Value = 1.983568429946899
$ws.Range("E1").Value = 0.7617396116256714
